# Updates the cryptos list sheet with refreshed price/volume data
# (mirrors the "Updated cryptos list ... with GitHub Actions" commit).
#
# Price (D) cells are plain decimal-looking strings that Excel would
# otherwise auto-convert to numbers (dropping trailing zeros etc.), so we
# briefly force a text number-format while assigning them, then restore
# the cell's original style so no formatting changes leak into the file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.400.38"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  +2.69%  "
$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.983.82"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  -0.09%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.94"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +1.96%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.24"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +3.37%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +0.89%  "
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.978.97"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("E10").Value = "  +2.93%  "
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.39"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  +11.53%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.450"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +0.10%  "
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000229"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  +3.90%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.73"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  +2.74%  "
$ws.Range("E15").Value = "  +0.06%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.476.66"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +1.12%  "
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.08"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  +1.38%  "
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.985.10"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  +1.26%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "59.420.78"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +2.79%  "
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "438.22"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +4.88%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.61"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +1.40%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.721"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +3.03%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.02"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -0.06%  "
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.24"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -1.37%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.01"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("E26").Value = "  -0.10%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.23"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +9.66%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("E30").Value = "  +2.76%  "
$ws.Range("E31").Value = "  +4.75%  "
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.81"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("E33").Value = "  +7.76%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0775"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  +10.41%  "
$ws.Range("E35").Value = "  +3.19%  "
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.986"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +3.55%  "
$ws.Range("E37").Value = "  +0.55%  "
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.68"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +1.08%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.68"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -2.65%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.76"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  +1.62%  "
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "400.94"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +4.05%  "
$ws.Range("E42").Value = "  +0.77%  "
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.735.32"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("E45").Value = "  +5.59%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.73"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +18.59%  "
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("E50").Value = "  +0.59%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.24"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +1.54%  "
